# zss447 - dateFormatDisplay testcase
# Populate sheet1 ("工作表1") with four rows that show a date value (41632)
# formatted in different ways, alongside a descriptive text label in column B.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Column A: date values + number formats -------------------------------
# Row 1 & 2 intentionally share the workbook's pre-existing "m/d/yyyy"-style
# date format (numFmtId 14) - write the plain value first so the existing
# style is reused/shared rather than being overwritten in place.
$ws.Range("A1").Value = 41632
$ws.Range("A2").Value = 41632
$ws.Range("A2").NumberFormat = "m/d/yy"

# Row 3 uses an explicit custom "m/d/yyyy" numeric format.
$ws.Range("A3").Value = 41632
$ws.Range("A3").NumberFormat = "m/d/yyyy"

# Row 4 uses a distinct "m/d/yy" style (kept separate from rows 1-2's style).
$ws.Range("A4").Value = 41632
$ws.Range("A4").NumberFormat = "m/d/yy;@"

# --- Column B: descriptive text labels -------------------------------------
$ws.Range("B3").Value = "<=m/d/yyyy"
$ws.Range("B2").Value = "<=yyyy/m/d, depends on locale"
$ws.Range("B1").Value = "<=No Format Configuration, depends on locale"
$ws.Range("B4").Value = "<=m/d/yy, doesn't support yet"

# --- Column widths -----------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 12.142857142857142
$ws.Columns.Item(2).ColumnWidth = 26

# --- Selection ------------------------------------------------------------
$null = $ws.Range("F8").Select()
